# Removed Test Case Inter-Dependency
$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# --- ProductLoanInput sheet (Sheet1) ---
# B1: productname -> new value
$wsInput.Range("B1").Value = "4200-RBI-EI-DB-DL-REC-FEE-RNI-FFC-SAR-FFROP-DAILY-1-CTRFD-MD-TR-1-ONTIME-PER-1st"

# B2: shortname -> now a text value "420d" instead of numeric 4200
$wsInput.Range("B2").Value = "420d"

# Selection / view state on input sheet moves to B1, and it's no longer the
# "tabSelected" sheet (ProductLoanOutput becomes selected instead).
$wsInput.Range("B1").Select()

# --- ProductLoanOutput sheet (Sheet2) ---
$wsOutput.Range("B1").Value = "4200-RBI-EI-DB-DL-REC-FEE-RNI-FFC-SAR-FFROP-DAILY-1-CTRFD-MD-TR-1-ONTIME-PER-1st"

# Make ProductLoanOutput the active sheet/tab and select B1 on it.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
